$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H29").Value = 1899.5
$ws.Range("I29").Value = 1899.5
$ws.Range("K29").Value = 5698.5
$ws.Range("M29").Value = -5417.5
$ws.Range("H32").Value = 3102.9614
$ws.Range("I32").Value = 3459.5334
$ws.Range("K32").Value = 3459.5334
$ws.Range("M32").Value = -3133.5334
$ws.Range("H33").Value = 705.46155
$ws.Range("I33").Value = 338.81818
$ws.Range("K33").Value = 338.81818
$ws.Range("M33").Value = -109.81818
$ws.Range("H38").Value = 1137.6154
$ws.Range("J38").Value = 12500
$ws.Range("L38").Value = 37500
$ws.Range("N38").Value = -38244
$ws.Range("H40").Value = 3172.2104
$ws.Range("J40").Value = 3065.1667
$ws.Range("L40").Value = 3065.1667
$ws.Range("N40").Value = -3415.1667
$ws.Range("H111").Value = 5999.25
$ws.Range("I111").Value = 5599
$ws.Range("J111").Value = 6399.5
$ws.Range("K111").Value = 16797
$ws.Range("L111").Value = 19198.5
$ws.Range("M111").Value = -13730
$ws.Range("N111").Value = -25332.5
$ws.Range("H116").Value = 10108392
$ws.Range("I116").Value = 18524912
$ws.Range("K116").Value = 18524912
$ws.Range("M116").Value = -18521470
$ws.Range("H125").Value = 3820.1667
$ws.Range("I125").Value = 2774.3333
$ws.Range("J125").Value = 4343.0835
$ws.Range("K125").Value = 24968.9997
$ws.Range("L125").Value = 39087.7515
$ws.Range("M125").Value = -22508.9997
$ws.Range("N125").Value = -44007.7515
$ws.Range("H135").Value = 3855.6875
$ws.Range("I135").Value = 4744.636
$ws.Range("J135").Value = 1900
$ws.Range("K135").Value = 42701.724
$ws.Range("L135").Value = 17100
$ws.Range("M135").Value = -40166.724
$ws.Range("N135").Value = -22170

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5993.0835
$ws.Range("I61").Value = 7259
$ws.Range("K61").Value = 7259
$ws.Range("M61").Value = -7047
$ws.Range("H74").Value = 3214.6667
$ws.Range("I74").Value = 1590.4667
$ws.Range("K74").Value = 1590.4667
$ws.Range("M74").Value = -716.4666999999999
$ws.Range("H77").Value = 3214.6667
$ws.Range("I77").Value = 1590.4667
$ws.Range("K77").Value = 7952.3335
$ws.Range("M77").Value = -3584.3335
$ws.Range("H122").Value = 522284.78
$ws.Range("I122").Value = 2683
$ws.Range("K122").Value = 8049
$ws.Range("M122").Value = -5599
$ws.Range("H132").Value = 3677.625
$ws.Range("I132").Value = 2137.125
$ws.Range("K132").Value = 6411.375
$ws.Range("M132").Value = -3881.375
$ws.Range("H136").Value = 5993.0835
$ws.Range("I136").Value = 7259
$ws.Range("K136").Value = 21777
$ws.Range("M136").Value = -19227

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 22138.059
$ws.Range("I99").Value = 33269.547
$ws.Range("J99").Value = 1730.3334
$ws.Range("K99").Value = 33269.547
$ws.Range("L99").Value = 1730.3334
$ws.Range("M99").Value = -31771.547
$ws.Range("N99").Value = -4726.3334
$ws.Range("H138").Value = 112500
$ws.Range("J138").Value = 112500
$ws.Range("L138").Value = 112500
$ws.Range("N138").Value = -122780
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").ClearContents()
$ws.Range("N139").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 971.7
$ws.Range("J22").Value = 985
$ws.Range("L22").Value = 985
$ws.Range("N22").Value = -1685
$ws.Range("H31").Value = 4908.7617
$ws.Range("I31").Value = 3835.077
$ws.Range("K31").Value = 3835.077
$ws.Range("M31").Value = -3540.077
$ws.Range("H34").Value = 4908.7617
$ws.Range("I34").Value = 3835.077
$ws.Range("K34").Value = 3835.077
$ws.Range("M34").Value = -3633.077
$ws.Range("H58").Value = 3170.7666
$ws.Range("I58").Value = 1884.4375
$ws.Range("J58").Value = 4640.857
$ws.Range("K58").Value = 1884.4375
$ws.Range("L58").Value = 4640.857
$ws.Range("M58").Value = -1681.4375
$ws.Range("N58").Value = -5046.857
$ws.Range("H122").Value = 1477.1666
$ws.Range("I122").Value = 1462.25
$ws.Range("K122").Value = 4386.75
$ws.Range("M122").Value = -1936.75
$ws.Range("H136").Value = 3170.7666
$ws.Range("I136").Value = 1884.4375
$ws.Range("J136").Value = 4640.857
$ws.Range("K136").Value = 5653.3125
$ws.Range("L136").Value = 13922.571
$ws.Range("M136").Value = -3103.3125
$ws.Range("N136").Value = -19022.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 9229.046
$ws.Range("I140").Value = 9826.950000000001
$ws.Range("K140").Value = 29480.85
$ws.Range("M140").Value = -24300.85

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 16279111
$ws.Range("I11").Value = 23583334
$ws.Range("J11").Value = 1670666.5
$ws.Range("K11").Value = 23583334
$ws.Range("L11").Value = 1670666.5
$ws.Range("M11").Value = -23583195
$ws.Range("N11").Value = -1670944.5
$ws.Range("H93").Value = 47777
$ws.Range("J93").Value = 47777
$ws.Range("L93").Value = 47777
$ws.Range("N93").Value = -51521
$ws.Range("H122").Value = 26613
$ws.Range("J122").Value = 33500
$ws.Range("L122").Value = 100500
$ws.Range("N122").Value = -105400

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 21975.682
$ws.Range("I7").Value = 33559.69
$ws.Range("J7").Value = 5243.222
$ws.Range("K7").Value = 33559.69
$ws.Range("L7").Value = 5243.222
$ws.Range("M7").Value = -33447.69
$ws.Range("N7").Value = -5467.222
$ws.Range("H22").Value = 4231.5713
$ws.Range("I22").Value = 5811
$ws.Range("J22").Value = 3599.8
$ws.Range("K22").Value = 5811
$ws.Range("L22").Value = 3599.8
$ws.Range("M22").Value = -5516
$ws.Range("N22").Value = -4189.8
$ws.Range("H27").Value = 4231.5713
$ws.Range("I27").Value = 5811
$ws.Range("J27").Value = 3599.8
$ws.Range("K27").Value = 5811
$ws.Range("L27").Value = 3599.8
$ws.Range("M27").Value = -5704
$ws.Range("N27").Value = -3813.8
$ws.Range("H40").Value = 67566.28999999999
$ws.Range("I40").Value = 102243.75
$ws.Range("K40").Value = 102243.75
$ws.Range("M40").Value = -102107.75
$ws.Range("H46").Value = 4116.8237
$ws.Range("J46").Value = 5590.273
$ws.Range("L46").Value = 5590.273
$ws.Range("N46").Value = -5966.273
$ws.Range("H55").Value = 930.2
$ws.Range("J55").Value = 1402.1428
$ws.Range("L55").Value = 1402.1428
$ws.Range("N55").Value = -1748.1428
$ws.Range("H100").Value = 9833.5
$ws.Range("I100").Value = 11200.2
$ws.Range("K100").Value = 11200.2
$ws.Range("M100").Value = -10659.2
$ws.Range("H122").Value = 5143.7856
$ws.Range("I122").Value = 4811.5
$ws.Range("K122").Value = 14434.5
$ws.Range("M122").Value = -11984.5
$ws.Range("H126").Value = 21975.682
$ws.Range("I126").Value = 33559.69
$ws.Range("J126").Value = 5243.222
$ws.Range("K126").Value = 100679.07
$ws.Range("L126").Value = 15729.666
$ws.Range("M126").Value = -98209.07000000001
$ws.Range("N126").Value = -20669.666
$ws.Range("H132").Value = 603210.8
$ws.Range("I132").Value = 885169.25
$ws.Range("K132").Value = 2655507.75
$ws.Range("M132").Value = -2652977.75
$ws.Range("H136").Value = 3787.4211
$ws.Range("J136").Value = 6696.4116
$ws.Range("L136").Value = 20089.2348
$ws.Range("N136").Value = -25189.2348
$ws.Range("H139").Value = 89499
$ws.Range("J139").Value = 89000
$ws.Range("L139").Value = 89000
$ws.Range("N139").Value = -99280

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20618.75
$ws.Range("I81").Value = 31140
$ws.Range("J81").Value = 3083.3333
$ws.Range("K81").Value = 62280
$ws.Range("L81").Value = 6166.6666
$ws.Range("M81").Value = -61219
$ws.Range("N81").Value = -8288.6666
$ws.Range("H84").Value = 20618.75
$ws.Range("I84").Value = 31140
$ws.Range("J84").Value = 3083.3333
$ws.Range("K84").Value = 311400
$ws.Range("L84").Value = 30833.333
$ws.Range("M84").Value = -306096
$ws.Range("N84").Value = -41441.333
$ws.Range("H100").Value = 32024.723
$ws.Range("I100").Value = 20342.076
$ws.Range("K100").Value = 40684.152
$ws.Range("M100").Value = -40143.152
$ws.Range("H107").Value = 37098.89
$ws.Range("I107").Value = 3770.1428
$ws.Range("K107").Value = 11310.4284
$ws.Range("M107").Value = -9390.428400000001
$ws.Range("H122").Value = 6340.6875
$ws.Range("I122").Value = 3299.8823
$ws.Range("J122").Value = 9786.933999999999
$ws.Range("K122").Value = 9899.6469
$ws.Range("L122").Value = 29360.802
$ws.Range("M122").Value = -7449.6469
$ws.Range("N122").Value = -34260.802
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
